$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 11 ("2021年") to the bottom of the data table, matching the
# existing table's layout/style (columns A:AQ).

# Copy the formatting of the existing year-label cell (A10) onto A11 so the
# new label keeps the same bold/centered/bordered style used by the other
# year rows in column A.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 4877.91
$ws.Range("C11").Value = 1105.47
$ws.Range("D11").Value = 260.8
# Column E is blank for this row (same as the other rows in the table).
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = 2075.79
$ws.Range("G11").Value = 4520.73
$ws.Range("H11").Value = 739.92
$ws.Range("I11").Value = 3565.27
$ws.Range("J11").Value = 321.69
$ws.Range("K11").Value = 97426.03
$ws.Range("L11").Value = 478.14
$ws.Range("M11").Value = 235.97
$ws.Range("N11").Value = 68.43000000000001
$ws.Range("O11").Value = 1448.59
$ws.Range("P11").Value = 5130.69
$ws.Range("Q11").Value = 152.01
$ws.Range("R11").Value = 217.95
$ws.Range("S11").Value = 1686.91
$ws.Range("T11").Value = 130.41
$ws.Range("U11").Value = 6404.46
$ws.Range("V11").Value = 3979.7
$ws.Range("W11").Value = 829.76
$ws.Range("X11").Value = 173.74
$ws.Range("Y11").Value = 1537.44
$ws.Range("Z11").Value = 7575.11
$ws.Range("AA11").Value = 509.09
$ws.Range("AB11").Value = 4851.13
$ws.Range("AC11").Value = 175.99
$ws.Range("AD11").Value = 1626.74
$ws.Range("AE11").Value = 1013.6
$ws.Range("AF11").Value = 15730.63
$ws.Range("AG11").Value = 5261.16
$ws.Range("AH11").Value = 989.8099999999999
$ws.Range("AI11").Value = 2549.94
$ws.Range("AJ11").Value = 230.81
$ws.Range("AK11").Value = 2802.39
$ws.Range("AL11").Value = 3611.11
$ws.Range("AM11").Value = 2420.67
$ws.Range("AN11").Value = 53.78
$ws.Range("AO11").Value = 1107.68
$ws.Range("AP11").Value = 6821.55
$ws.Range("AQ11").Value = 152.23
